$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case TC-113 - fill column by column (A, then B..D headers/data,
# then E, F, G) so new shared-string entries are interned in the same
# order as the source workbook.
$ws.Range("A7").Value = "TC-113"
$ws.Range("A8").Value = "TC-113"

$ws.Range("B7").Value = "username"
$ws.Range("B8").Value = "reyaz0806"

$ws.Range("C7").Value = "password"
$ws.Range("C8").Value = "reyaz123"

$ws.Range("D7").Value = "expected Title"
$ws.Range("D8").Value = "Adactin.com - Search Hotel"

$ws.Range("E7").Value = "location"
$ws.Range("E8").Value = "Sydney"

$ws.Range("F7").Value = "hotel"
$ws.Range("F8").Value = "Hotel Creek"

$ws.Range("G7").Value = "Check In Date"

# Date typed with a leading apostrophe (stored as text w/ quote-prefix) but
# displayed with a date number format, matching the source workbook.
$ws.Range("G8").Value = "'27/09/2024"
$ws.Range("G8").NumberFormat = "mm-dd-yy"

# Match the new best-fit column widths for the two added columns.
$ws.Columns.Item(6).ColumnWidth = 9.833333333333334
$ws.Columns.Item(7).ColumnWidth = 11.5

# Move the active selection to G7, as in the edited workbook.
$ws.Range("G7").Select() | Out-Null
